# Generate Report for Archive
#
# The Overview / zh-cn / de-de sheets list one row per source file
# (709dd644..., 1379beff..., b7a24cfc..., d2a702dc...). The b7a24cfc and
# d2a702dc rows have moved into "In Translation" and are now sorted ahead of
# the still-pending 1379beff row, which drops to the bottom. This script
# rewrites the row 3/4/5 cell contents on all three sheets to reflect the
# new ordering/status, matching the updated shared-strings table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 3: now b7a24cfc, In Translation, handed off at 09:41:43
$ws.Range("A3").Value = "b7a24cfc-1868-49ee-8e67-8ceaf538fd80.md"
$ws.Range("B3").Value = "In Translation"
$ws.Range("C3").Value = "In Translation"
$ws.Range("D3").Value = "2016-03-24 09:41:43"

# Row 4: now d2a702dc, In Translation, handed off at 09:41:43
$ws.Range("A4").Value = "d2a702dc-6c9a-4ead-be34-f2a38ed14158.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "2016-03-24 09:41:43"

# Row 5: now 1379beff, still Ready for handoff at 09:39:42
$ws.Range("A5").Value = "1379beff-452b-4dcd-a8bd-6b4d0df301da.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "2016-03-24 09:39:42"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 3: now b7a24cfc
$ws.Range("A3").Value = "b7a24cfc-1868-49ee-8e67-8ceaf538fd80.md"
$ws.Range("C3").Value = "In Translation"
$ws.Range("D3").Value = "b7a24cfc-1868-49ee-8e67-8ceaf538fd80.9a90ab762dfa2221dd4ae971866e00a32018abe9.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-24 09:41:39"

# Row 4: now d2a702dc
$ws.Range("A4").Value = "d2a702dc-6c9a-4ead-be34-f2a38ed14158.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "d2a702dc-6c9a-4ead-be34-f2a38ed14158.b5e2b64e353d982114dc509e0c5cf7a4ce07bb7e.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-24 09:41:39"

# Row 5: now 1379beff
$ws.Range("A5").Value = "1379beff-452b-4dcd-a8bd-6b4d0df301da.md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "1379beff-452b-4dcd-a8bd-6b4d0df301da.d02ce7e7f9a3c8b6a4c080d188d3343c32b883fa.zh-cn.xlf"
$ws.Range("E5").Value = "2016-03-24 09:39:38"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 3: now b7a24cfc
$ws.Range("A3").Value = "b7a24cfc-1868-49ee-8e67-8ceaf538fd80.md"
$ws.Range("C3").Value = "In Translation"
$ws.Range("D3").Value = "b7a24cfc-1868-49ee-8e67-8ceaf538fd80.9a90ab762dfa2221dd4ae971866e00a32018abe9.de-de.xlf"
$ws.Range("E3").Value = "2016-03-24 09:41:43"

# Row 4: now d2a702dc
$ws.Range("A4").Value = "d2a702dc-6c9a-4ead-be34-f2a38ed14158.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "d2a702dc-6c9a-4ead-be34-f2a38ed14158.b5e2b64e353d982114dc509e0c5cf7a4ce07bb7e.de-de.xlf"
$ws.Range("E4").Value = "2016-03-24 09:41:43"

# Row 5: now 1379beff
$ws.Range("A5").Value = "1379beff-452b-4dcd-a8bd-6b4d0df301da.md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "1379beff-452b-4dcd-a8bd-6b4d0df301da.d02ce7e7f9a3c8b6a4c080d188d3343c32b883fa.de-de.xlf"
$ws.Range("E5").Value = "2016-03-24 09:39:42"
